$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5.27
$ws.Range("D3").Value = 0.22
$ws.Range("D4").Value = 0.93
$ws.Range("D6").Value = 0.2
$ws.Range("D7").Value = 0.1
$ws.Range("D8").Value = 0.19
$ws.Range("D9").Value = 0.15
$ws.Range("D10").Value = 0.08
$ws.Range("D11").Value = 0.23
$ws.Range("D12").Value = 4.51
$ws.Range("D13").Value = 0.79
$ws.Range("D14").Value = 13.2
$ws.Range("D15").Value = 2.63
$ws.Range("D16").Value = 0.23
$ws.Range("D17").Value = 7.21
$ws.Range("D18").Value = 0.13
$ws.Range("D19").Value = 0.36
$ws.Range("D21").Value = 0.1
$ws.Range("D22").Value = 0.9399999999999999
$ws.Range("D23").Value = 0.71
